$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.201.28'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.602.33'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.39'
$ws.Range("E6").Value = '  +0.54%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.07'
$ws.Range("E8").Value = '  +4.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3616'
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("E10").Value = '  +0.19%  '
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08118'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.75'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.584'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.420'
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").Value = '1.599.79'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.21'
$ws.Range("E18").Value = '  +2.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06883'
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.550'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.99'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").Value = '23.174.61'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("E25").Value = '  +2.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.983'
$ws.Range("E26").Value = '  +9.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.25'
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.47'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.247'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.02'
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.772'
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").Value = '1.779.28'
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07510'
$ws.Range("E35").Value = '  -2.46%  '
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02715'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2504'
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.08808'
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.079'
$ws.Range("E40").Value = '  -3.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7109'
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.360'
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.63'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6525'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.017'
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.18'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07971'
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  -3.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.215'
$ws.Range("E51").Value = '  +1.62%  '
